# Applies the updated crypto price/volume figures (and two name/link
# row swaps) from the Fri Jan 19 07:15:14 UTC 2024 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to stay text (matches the source data, which stores
    # prices/percentages as strings, not numbers) without leaving a
    # residual number-format behind on the cell once written.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Row 2
Set-TextValue $ws.Range("D2") "41.498.12"
Set-TextValue $ws.Range("E2") "  -3.28%  "

# Row 3
Set-TextValue $ws.Range("D3") "2.470.80"
Set-TextValue $ws.Range("E3") "  -2.52%  "

# Row 4
Set-TextValue $ws.Range("D4") "1.00"
Set-TextValue $ws.Range("E4") "  +0.04%  "

# Row 5
Set-TextValue $ws.Range("D5") "312.00"
Set-TextValue $ws.Range("E5") "  -0.10%  "

# Row 6
Set-TextValue $ws.Range("D6") "94.89"
Set-TextValue $ws.Range("E6") "  -6.34%  "

# Row 7
Set-TextValue $ws.Range("E7") "  -2.73%  "

# Row 8
Set-TextValue $ws.Range("E8") "  +0.04%  "

# Row 9
Set-TextValue $ws.Range("E9") "  -4.62%  "

# Row 10
Set-TextValue $ws.Range("D10") "33.60"
Set-TextValue $ws.Range("E10") "  -6.46%  "

# Row 11
Set-TextValue $ws.Range("E11") "  -3.14%  "

# Row 12
Set-TextValue $ws.Range("E12") "  -1.03%  "

# Row 13
Set-TextValue $ws.Range("D13") "7.04"
Set-TextValue $ws.Range("E13") "  -4.30%  "

# Row 14
Set-TextValue $ws.Range("D14") "2.851.44"
Set-TextValue $ws.Range("E14") "  -2.47%  "

# Row 15
Set-TextValue $ws.Range("D15") "2.446.05"
Set-TextValue $ws.Range("E15") "  -3.70%  "

# Row 16
Set-TextValue $ws.Range("D16") "15.00"
Set-TextValue $ws.Range("E16") "  -2.84%  "

# Row 17
Set-TextValue $ws.Range("E17") "  -3.72%  "

# Row 18
Set-TextValue $ws.Range("D18") "41.504.16"
Set-TextValue $ws.Range("E18") "  -3.17%  "

# Row 19
Set-TextValue $ws.Range("D19") "6.35"
Set-TextValue $ws.Range("E19") "  -5.09%  "

# Row 20
Set-TextValue $ws.Range("D20") "0.0₃0924"
Set-TextValue $ws.Range("E20") "  -3.25%  "

# Row 21
Set-TextValue $ws.Range("D21") "11.27"
Set-TextValue $ws.Range("E21") "  -9.27%  "

# Row 22
Set-TextValue $ws.Range("E22") "  -1.83%  "

# Row 23
Set-TextValue $ws.Range("D23") "237.50"
Set-TextValue $ws.Range("E23") "  -2.85%  "

# Row 24
Set-TextValue $ws.Range("D24") "2.76"
Set-TextValue $ws.Range("E24") "  -4.47%  "

# Row 25
Set-TextValue $ws.Range("B25") "ImmutableX"
Set-TextValue $ws.Range("C25") "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D25") "1.91"
Set-TextValue $ws.Range("E25") "  -6.38%  "

# Row 26
Set-TextValue $ws.Range("B26") "Dai"
Set-TextValue $ws.Range("C26") "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Range("D26") "1.00"
Set-TextValue $ws.Range("E26") "  -0.10%  "

# Row 27
Set-TextValue $ws.Range("D27") "24.20"
Set-TextValue $ws.Range("E27") "  -5.57%  "

# Row 28
Set-TextValue $ws.Range("D28") "2.22"
Set-TextValue $ws.Range("E28") "  -5.39%  "

# Row 29
Set-TextValue $ws.Range("D29") "9.66"
Set-TextValue $ws.Range("E29") "  -5.72%  "

# Row 30
Set-TextValue $ws.Range("D30") "36.74"
Set-TextValue $ws.Range("E30") "  -5.52%  "

# Row 31
Set-TextValue $ws.Range("D31") "152.18"
Set-TextValue $ws.Range("E31") "  -5.58%  "

# Row 32
Set-TextValue $ws.Range("D32") "5.50"
Set-TextValue $ws.Range("E32") "  -6.07%  "

# Row 33
Set-TextValue $ws.Range("D33") "2.66"
Set-TextValue $ws.Range("E33") "  -4.77%  "

# Row 34
Set-TextValue $ws.Range("E34") "  -2.25%  "

# Row 35
Set-TextValue $ws.Range("D35") "0.0748"
Set-TextValue $ws.Range("E35") "  -5.56%  "

# Row 36
Set-TextValue $ws.Range("D36") "3.05"
Set-TextValue $ws.Range("E36") "  -2.31%  "

# Row 37
Set-TextValue $ws.Range("B37") "Celestia"
Set-TextValue $ws.Range("C37") "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextValue $ws.Range("D37") "17.10"
Set-TextValue $ws.Range("E37") "  -7.05%  "

# Row 38
Set-TextValue $ws.Range("B38") "ARBITRUM"
Set-TextValue $ws.Range("C38") "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D38") "1.89"
Set-TextValue $ws.Range("E38") "  -4.47%  "

# Row 39
Set-TextValue $ws.Range("E39") "  -2.98%  "

# Row 40
Set-TextValue $ws.Range("B40") "Kaspa"
Set-TextValue $ws.Range("C40") "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D40") "0.103"
Set-TextValue $ws.Range("E40") "  -8.09%  "

# Row 41
Set-TextValue $ws.Range("B41") "RenderToken"
Set-TextValue $ws.Range("C41") "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D41") "4.27"
Set-TextValue $ws.Range("E41") "  +1.84%  "

# Row 42
Set-TextValue $ws.Range("E42") "  +0.18%  "

# Row 43
Set-TextValue $ws.Range("D43") "20.04"
Set-TextValue $ws.Range("E43") "  -9.06%  "

# Row 44
Set-TextValue $ws.Range("D44") "1.989.45"
Set-TextValue $ws.Range("E44") "  -0.57%  "

# Row 45
Set-TextValue $ws.Range("E45") "  -4.41%  "

# Row 46
Set-TextValue $ws.Range("E46") "  -8.94%  "

# Row 47
Set-TextValue $ws.Range("D47") "8.78"
Set-TextValue $ws.Range("E47") "  -4.61%  "

# Row 48
Set-TextValue $ws.Range("D48") "2.716.17"
Set-TextValue $ws.Range("E48") "  -2.09%  "

# Row 49
Set-TextValue $ws.Range("D49") "69.96"
Set-TextValue $ws.Range("E49") "  -3.72%  "

# Row 50
Set-TextValue $ws.Range("D50") "96.77"
Set-TextValue $ws.Range("E50") "  -4.73%  "

# Row 51
Set-TextValue $ws.Range("D51") "75.12"
Set-TextValue $ws.Range("E51") "  -5.81%  "
